$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare formatting for the 3 new rows (47-49) by copying column-A style from row 46
$ws.Range("A46").Copy() | Out-Null
$ws.Range("A47:A49").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(2, 1).Value = 'Daxak'
$ws.Cells.Item(2, 2).Value = 227.319653937568
$ws.Cells.Item(2, 3).Value = 3
$ws.Cells.Item(2, 4).Value = 15
$ws.Cells.Item(2, 5).Value = 15.15464359583787

$ws.Cells.Item(3, 1).Value = 'MagicaL'
$ws.Cells.Item(3, 2).Value = 393.6062553980921
$ws.Cells.Item(3, 3).Value = 16
$ws.Cells.Item(3, 4).Value = 29
$ws.Cells.Item(3, 5).Value = 13.57262949648593

$ws.Cells.Item(4, 1).Value = 'Chessie'
$ws.Cells.Item(4, 2).Value = 143.815580533439
$ws.Cells.Item(4, 3).Value = 8
$ws.Cells.Item(4, 4).Value = 11
$ws.Cells.Item(4, 5).Value = 13.07414368485809

$ws.Cells.Item(5, 1).Value = 'Abed'
$ws.Cells.Item(5, 2).Value = 622.6148198856745
$ws.Cells.Item(5, 3).Value = 36
$ws.Cells.Item(5, 4).Value = 50
$ws.Cells.Item(5, 5).Value = 12.45229639771349

$ws.Cells.Item(6, 1).Value = 201594424
$ws.Cells.Item(6, 2).Value = 210.2644480396311
$ws.Cells.Item(6, 3).Value = 8
$ws.Cells.Item(6, 4).Value = 17
$ws.Cells.Item(6, 5).Value = 12.36849694350771

$ws.Cells.Item(7, 1).Value = 'Nisha'
$ws.Cells.Item(7, 2).Value = 194.9083048153606
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(7, 4).Value = 16
$ws.Cells.Item(7, 5).Value = 12.18176905096004

$ws.Cells.Item(8, 1).Value = 'J E I M A R I '
$ws.Cells.Item(8, 2).Value = 84.33886859741604
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 4).Value = 7
$ws.Cells.Item(8, 5).Value = 12.04840979963086

$ws.Cells.Item(9, 1).Value = '小可'
$ws.Cells.Item(9, 2).Value = 503.5703066322106
$ws.Cells.Item(9, 3).Value = 28
$ws.Cells.Item(9, 4).Value = 42
$ws.Cells.Item(9, 5).Value = 11.98976920552882

$ws.Cells.Item(10, 1).Value = 'SumaiL'
$ws.Cells.Item(10, 2).Value = 795.1177358204511
$ws.Cells.Item(10, 3).Value = 36
$ws.Cells.Item(10, 4).Value = 67
$ws.Cells.Item(10, 5).Value = 11.86742889284255

$ws.Cells.Item(11, 1).Value = 'Limmp'
$ws.Cells.Item(11, 2).Value = 236.7331774912908
$ws.Cells.Item(11, 3).Value = 11
$ws.Cells.Item(11, 4).Value = 20
$ws.Cells.Item(11, 5).Value = 11.83665887456454

$ws.Cells.Item(12, 1).Value = 'fn'
$ws.Cells.Item(12, 2).Value = 141.4301012369334
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(12, 4).Value = 12
$ws.Cells.Item(12, 5).Value = 11.78584176974445

$ws.Cells.Item(13, 1).Value = '430'
$ws.Cells.Item(13, 2).Value = 223.7306337688624
$ws.Cells.Item(13, 3).Value = 3
$ws.Cells.Item(13, 4).Value = 19
$ws.Cells.Item(13, 5).Value = 11.77529651415065

$ws.Cells.Item(14, 1).Value = 'AhJit'
$ws.Cells.Item(14, 2).Value = 278.1358160058417
$ws.Cells.Item(14, 3).Value = 13
$ws.Cells.Item(14, 4).Value = 24
$ws.Cells.Item(14, 5).Value = 11.58899233357674

$ws.Cells.Item(15, 1).Value = 'Setsu'
$ws.Cells.Item(15, 2).Value = 182.6941839397884
$ws.Cells.Item(15, 3).Value = 6
$ws.Cells.Item(15, 4).Value = 16
$ws.Cells.Item(15, 5).Value = 11.41838649623678

$ws.Cells.Item(16, 1).Value = 'ARMEL'
$ws.Cells.Item(16, 2).Value = 262.5496595994638
$ws.Cells.Item(16, 3).Value = 25
$ws.Cells.Item(16, 4).Value = 23
$ws.Cells.Item(16, 5).Value = 11.41520259128103

$ws.Cells.Item(17, 1).Value = 'Resolut1on'
$ws.Cells.Item(17, 2).Value = 125.5188420054152
$ws.Cells.Item(17, 3).Value = 5
$ws.Cells.Item(17, 4).Value = 11
$ws.Cells.Item(17, 5).Value = 11.41080381867411

$ws.Cells.Item(18, 1).Value = 'ASD'
$ws.Cells.Item(18, 2).Value = 403.5494247930765
$ws.Cells.Item(18, 3).Value = 31
$ws.Cells.Item(18, 4).Value = 36
$ws.Cells.Item(18, 5).Value = 11.20970624425212

$ws.Cells.Item(19, 1).Value = 'No[o]ne-'
$ws.Cells.Item(19, 2).Value = 851.056595955431
$ws.Cells.Item(19, 3).Value = 36
$ws.Cells.Item(19, 4).Value = 76
$ws.Cells.Item(19, 5).Value = 11.19811310467672

$ws.Cells.Item(20, 1).Value = 'k`wonderkid'
$ws.Cells.Item(20, 2).Value = 264.446722416145
$ws.Cells.Item(20, 3).Value = 13
$ws.Cells.Item(20, 4).Value = 24
$ws.Cells.Item(20, 5).Value = 11.01861343400604

$ws.Cells.Item(21, 1).Value = '_Mikoto_'
$ws.Cells.Item(21, 2).Value = 294.7336774507343
$ws.Cells.Item(21, 3).Value = 9
$ws.Cells.Item(21, 4).Value = 27
$ws.Cells.Item(21, 5).Value = 10.91606212780498

$ws.Cells.Item(22, 1).Value = 'Xm'
$ws.Cells.Item(22, 2).Value = 139.984847976795
$ws.Cells.Item(22, 3).Value = 8
$ws.Cells.Item(22, 4).Value = 13
$ws.Cells.Item(22, 5).Value = 10.76806522898423

$ws.Cells.Item(23, 1).Value = 'Ori'
$ws.Cells.Item(23, 2).Value = 750.8454423568251
$ws.Cells.Item(23, 3).Value = 39
$ws.Cells.Item(23, 4).Value = 70
$ws.Cells.Item(23, 5).Value = 10.72636346224036

$ws.Cells.Item(24, 1).Value = 'Chris Luck '
$ws.Cells.Item(24, 2).Value = 148.0914112609426
$ws.Cells.Item(24, 3).Value = 3
$ws.Cells.Item(24, 4).Value = 14
$ws.Cells.Item(24, 5).Value = 10.57795794721018

$ws.Cells.Item(25, 1).Value = 'Fata'
$ws.Cells.Item(25, 2).Value = 807.9741470712329
$ws.Cells.Item(25, 3).Value = 22
$ws.Cells.Item(25, 4).Value = 77
$ws.Cells.Item(25, 5).Value = 10.49317074118484

$ws.Cells.Item(26, 1).Value = 'Bryle'
$ws.Cells.Item(26, 2).Value = 261.8596368979141
$ws.Cells.Item(26, 3).Value = 33
$ws.Cells.Item(26, 4).Value = 25
$ws.Cells.Item(26, 5).Value = 10.47438547591656

$ws.Cells.Item(27, 1).Value = 'MidOne'
$ws.Cells.Item(27, 2).Value = 707.9299802840181
$ws.Cells.Item(27, 3).Value = 36
$ws.Cells.Item(27, 4).Value = 68
$ws.Cells.Item(27, 5).Value = 10.41073500417674

$ws.Cells.Item(28, 1).Value = 'qojqva'
$ws.Cells.Item(28, 2).Value = 380.290798259293
$ws.Cells.Item(28, 3).Value = 28
$ws.Cells.Item(28, 4).Value = 38
$ws.Cells.Item(28, 5).Value = 10.00765258577087

$ws.Cells.Item(29, 1).Value = 'lover'
$ws.Cells.Item(29, 2).Value = 148.6794197813322
$ws.Cells.Item(29, 3).Value = 25
$ws.Cells.Item(29, 4).Value = 15
$ws.Cells.Item(29, 5).Value = 9.911961318755482

$ws.Cells.Item(30, 1).Value = 'Afoninje'
$ws.Cells.Item(30, 2).Value = 322.0918559628618
$ws.Cells.Item(30, 3).Value = 11
$ws.Cells.Item(30, 4).Value = 33
$ws.Cells.Item(30, 5).Value = 9.760359271601873

$ws.Cells.Item(31, 1).Value = 'Moonn'
$ws.Cells.Item(31, 2).Value = 214.5299565160425
$ws.Cells.Item(31, 3).Value = 11
$ws.Cells.Item(31, 4).Value = 22
$ws.Cells.Item(31, 5).Value = 9.751361659820116

$ws.Cells.Item(32, 1).Value = 'G'
$ws.Cells.Item(32, 2).Value = 106.9541087173019
$ws.Cells.Item(32, 3).Value = 3
$ws.Cells.Item(32, 4).Value = 11
$ws.Cells.Item(32, 5).Value = 9.723100792481988

$ws.Cells.Item(33, 1).Value = 'CCnC'
$ws.Cells.Item(33, 2).Value = 211.169981931889
$ws.Cells.Item(33, 3).Value = 8
$ws.Cells.Item(33, 4).Value = 22
$ws.Cells.Item(33, 5).Value = 9.598635542358592

$ws.Cells.Item(34, 1).Value = '4dr <3 Déia'
$ws.Cells.Item(34, 2).Value = 66.31480507374468
$ws.Cells.Item(34, 3).Value = 3
$ws.Cells.Item(34, 4).Value = 7
$ws.Cells.Item(34, 5).Value = 9.473543581963526

$ws.Cells.Item(35, 1).Value = 'Dendi'
$ws.Cells.Item(35, 2).Value = 84.47633164830921
$ws.Cells.Item(35, 3).Value = 8
$ws.Cells.Item(35, 4).Value = 9
$ws.Cells.Item(35, 5).Value = 9.386259072034356

$ws.Cells.Item(36, 1).Value = '| Draken-_-'
$ws.Cells.Item(36, 2).Value = 74.2565603871312
$ws.Cells.Item(36, 3).Value = 20
$ws.Cells.Item(36, 4).Value = 8
$ws.Cells.Item(36, 5).Value = 9.2820700483914

$ws.Cells.Item(37, 1).Value = '一'
$ws.Cells.Item(37, 2).Value = 354.0692561586003
$ws.Cells.Item(37, 3).Value = 14
$ws.Cells.Item(37, 4).Value = 39
$ws.Cells.Item(37, 5).Value = 9.078698875861546

$ws.Cells.Item(38, 1).Value = 'w33'
$ws.Cells.Item(38, 2).Value = 299.8261416149526
$ws.Cells.Item(38, 3).Value = 36
$ws.Cells.Item(38, 4).Value = 34
$ws.Cells.Item(38, 5).Value = 8.818415929851547

$ws.Cells.Item(39, 1).Value = 'Topson'
$ws.Cells.Item(39, 2).Value = 277.0178318552783
$ws.Cells.Item(39, 3).Value = 9
$ws.Cells.Item(39, 4).Value = 32
$ws.Cells.Item(39, 5).Value = 8.656807245477447

$ws.Cells.Item(40, 1).Value = 'MATUMBAMAN'
$ws.Cells.Item(40, 2).Value = 381.1183383920733
$ws.Cells.Item(40, 3).Value = 31
$ws.Cells.Item(40, 4).Value = 45
$ws.Cells.Item(40, 5).Value = 8.46929640871274

$ws.Cells.Item(41, 1).Value = 'YawaR'
$ws.Cells.Item(41, 2).Value = 135.0902426416654
$ws.Cells.Item(41, 3).Value = 28
$ws.Cells.Item(41, 4).Value = 16
$ws.Cells.Item(41, 5).Value = 8.44314016510409

$ws.Cells.Item(42, 1).Value = 'p4pita'
$ws.Cells.Item(42, 2).Value = 109.7528918787175
$ws.Cells.Item(42, 3).Value = 8
$ws.Cells.Item(42, 4).Value = 13
$ws.Cells.Item(42, 5).Value = 8.442530144516732

$ws.Cells.Item(43, 1).Value = 'Ryoya'
$ws.Cells.Item(43, 2).Value = 210.3753633477911
$ws.Cells.Item(43, 3).Value = 9
$ws.Cells.Item(43, 4).Value = 25
$ws.Cells.Item(43, 5).Value = 8.415014533911645

$ws.Cells.Item(44, 1).Value = 'iAnnihilate'
$ws.Cells.Item(44, 2).Value = 49.57841213490627
$ws.Cells.Item(44, 3).Value = 8
$ws.Cells.Item(44, 4).Value = 6
$ws.Cells.Item(44, 5).Value = 8.263068689151046

$ws.Cells.Item(45, 1).Value = 'Timado'
$ws.Cells.Item(45, 2).Value = 63.05210138462046
$ws.Cells.Item(45, 3).Value = 8
$ws.Cells.Item(45, 4).Value = 8
$ws.Cells.Item(45, 5).Value = 7.881512673077557

$ws.Cells.Item(46, 1).Value = 'kodos-'
$ws.Cells.Item(46, 2).Value = 38.23544579454546
$ws.Cells.Item(46, 3).Value = 3
$ws.Cells.Item(46, 4).Value = 5
$ws.Cells.Item(46, 5).Value = 7.647089158909091

$ws.Cells.Item(47, 1).Value = 'hFn k3'
$ws.Cells.Item(47, 2).Value = 68.40668907883938
$ws.Cells.Item(47, 3).Value = 20
$ws.Cells.Item(47, 4).Value = 9
$ws.Cells.Item(47, 5).Value = 7.600743230982153

$ws.Cells.Item(48, 1).Value = 'Ceyler'
$ws.Cells.Item(48, 2).Value = 55.06828333169397
$ws.Cells.Item(48, 3).Value = 20
$ws.Cells.Item(48, 4).Value = 8
$ws.Cells.Item(48, 5).Value = 6.883535416461746

$ws.Cells.Item(49, 1).Value = 'Faker-'
$ws.Cells.Item(49, 2).Value = 27.03325952246283
$ws.Cells.Item(49, 3).Value = 3
$ws.Cells.Item(49, 4).Value = 6
$ws.Cells.Item(49, 5).Value = 4.505543253743805
